$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Clear obsolete columns (G/H entirely; D/E for rows 19-40) ---
$ws.Range("G15:H40").Clear() | Out-Null
$ws.Range("D19:E40").Clear() | Out-Null

# --- Re-apply existing named cell formats (copy format from an already-styled
#     cell so no NEW font/xf entries get created) ---
$ws.Range("A9").Copy() | Out-Null
foreach ($t in @("A10","A11","A12")) {
    $ws.Range($t).PasteSpecial(-4122) | Out-Null
}
$ws.Range("B10").Copy() | Out-Null
foreach ($t in @("A16","D16","A17","D17","A18","D18","A19","A20","A21","A23","A25","A26","A28","A29","A30","A32","A33","A36","A38","A39")) {
    $ws.Range($t).PasteSpecial(-4122) | Out-Null
}
$ws.Range("C10").Copy() | Out-Null
foreach ($t in @("A22","A24","A27","A31","A34","A35","A37","A40")) {
    $ws.Range($t).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Set cell values/text ---
$ws.Range("A10").Value2 = "No."
$ws.Range("B10").Value2 = 20
$ws.Range("C10").Value2 = 8
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 28
$ws.Range("A11").Value2 = "Marking"
$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1
$ws.Range("D11").Value2 = 0
$ws.Range("A12").Value2 = "Total"
$ws.Range("B12").Value2 = 80
$ws.Range("C12").Value2 = -8
$ws.Range("E12").Value2 = "72/112"
$ws.Range("A15").Value2 = "Student Ans"
$ws.Range("B15").Value2 = "Correct Ans"
$ws.Range("D15").Value2 = "Student Ans"
$ws.Range("E15").Value2 = "Correct Ans"
$ws.Range("A16").Value2 = "Option A"
$ws.Range("B16").Value2 = "Option A"
$ws.Range("D16").Value2 = "Option A"
$ws.Range("E16").Value2 = "Option A"
$ws.Range("A17").Value2 = "Option D"
$ws.Range("B17").Value2 = "Option D"
$ws.Range("D17").Value2 = "Option C"
$ws.Range("E17").Value2 = "Option C"
$ws.Range("A18").Value2 = "Option B"
$ws.Range("B18").Value2 = "Option B"
$ws.Range("D18").Value2 = "Option D"
$ws.Range("E18").Value2 = "Option D"
$ws.Range("A19").Value2 = "Option C"
$ws.Range("B19").Value2 = "Option C"
$ws.Range("A20").Value2 = "Option B"
$ws.Range("B20").Value2 = "Option B"
$ws.Range("A21").Value2 = "Option C"
$ws.Range("B21").Value2 = "Option C"
$ws.Range("A22").Value2 = "Option A"
$ws.Range("B22").Value2 = "Option D"
$ws.Range("A23").Value2 = "Option D"
$ws.Range("B23").Value2 = "Option D"
$ws.Range("A24").Value2 = "Option B"
$ws.Range("B24").Value2 = "Option A"
$ws.Range("A25").Value2 = "Option A"
$ws.Range("B25").Value2 = "Option A"
$ws.Range("A26").Value2 = "Option C"
$ws.Range("B26").Value2 = "Option C"
$ws.Range("A27").Value2 = "Option C"
$ws.Range("B27").Value2 = "Option A"
$ws.Range("A28").Value2 = "Option D"
$ws.Range("B28").Value2 = "Option D"
$ws.Range("A29").Value2 = "Option D"
$ws.Range("B29").Value2 = "Option D"
$ws.Range("A30").Value2 = "Option B"
$ws.Range("B30").Value2 = "Option B"
$ws.Range("A31").Value2 = "Option B"
$ws.Range("B31").Value2 = "Option D"
$ws.Range("A32").Value2 = "Option C"
$ws.Range("B32").Value2 = "Option C"
$ws.Range("A33").Value2 = "Option D"
$ws.Range("B33").Value2 = "Option D"
$ws.Range("A34").Value2 = "Option D"
$ws.Range("B34").Value2 = "Option B"
$ws.Range("A35").Value2 = "Option B"
$ws.Range("B35").Value2 = "Option D"
$ws.Range("A36").Value2 = "Option A"
$ws.Range("B36").Value2 = "Option A"
$ws.Range("A37").Value2 = "Option B"
$ws.Range("B37").Value2 = "Option A"
$ws.Range("A38").Value2 = "Option A"
$ws.Range("B38").Value2 = "Option A"
$ws.Range("A39").Value2 = "Option D"
$ws.Range("B39").Value2 = "Option D"
$ws.Range("A40").Value2 = "Option B"
$ws.Range("B40").Value2 = "Option D"
